$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student's username and email (e.g. account renamed with a "4" suffix)
$ws.Range("B2").Value = "TrianNurizkillah4"
$ws.Range("C2").Value = "triannurizkillah4@gmail.com"

# Update the active selection to reflect where the user last clicked (C2)
$ws.Range("C2").Select()
